$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B11: change its text value from "R40" to "1" (still stored as text, not a number)
$ws.Range("B11").Value = "'1"

